$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to be inserted right after the header row (row 1),
# pushing existing data rows down by 7.
$topRows = @(
    @(0.3322033213479865, 0.1210244501453987, -1.57638051243603),
    @(0.6352223770569996, -0.8322304098058368, -1.049409866333014),
    @(0.1475016089714127, -0.1921065187431543, -0.1445163721927848),
    @(-1.213809009395561, -0.2346711329983395, 1.28885372863315),
    @(-0.7734762763008961, -0.09618946024911929, 2.230298755737722),
    @(-0.814207781389884, 0.3517244944944581, 1.427822031950592),
    @(-0.4451152733739872, 0.1002355693144553, 1.227771341498127)
)

# Insert 7 new blank rows right after the header (before old row 2).
$insertRange = $ws.Range("A2:C8")
$insertRange.EntireRow.Insert()
# The inserted rows pick up formatting from the row above (the header);
# clear it so the new data rows stay unstyled like the rest of the data.
$ws.Range("A2:C8").ClearFormats()

for ($i = 0; $i -lt $topRows.Count; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $topRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $topRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $topRows[$i][2]
}

# New rows appended at the bottom (rows 29-31 after the insert above).
$bottomRows = @(
    @(-2.691276774793723, 7.824305781253101, -2.586581079487865),
    @(1.383527442585964, 2.926498572838676, -5.056850963437613),
    @(3.296718087898288, -3.505119464119212, -4.962326313638459)
)

for ($i = 0; $i -lt $bottomRows.Count; $i++) {
    $r = 29 + $i
    $ws.Cells.Item($r, 1).Value = $bottomRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $bottomRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $bottomRows[$i][2]
}
